# Refactor codes (create 2 Services)
# The only functional content change in the sheet is cell A2, which described
# "trueExploreMode" as a string parameter; it is now documented as a bool.
# We also move the active selection to C17 to mirror the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "         trueExploreMode: bool,"

$ws.Range("C17").Select()
